$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.613.24'
$ws.Range("E2").Value = '  +0.81%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.851.47'
$ws.Range("E3").Value = '  +0.34%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '264.67'
$ws.Range("E5").Value = '  +1.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5242'
$ws.Range("E7").Value = '  -0.30%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3251'
$ws.Range("E8").Value = '  +0.72%  '
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.97'
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7828'
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07792'
$ws.Range("E12").Value = '  +0.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.856.26'
$ws.Range("E13").Value = '  +1.01%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.65'
$ws.Range("E14").Value = '  -0.36%  '
$ws.Range("E15").Value = '  +0.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  -0.13%  '
$ws.Range("E17").Value = '  -0.89%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007979'
$ws.Range("E18").Value = '  +0.99%  '
$ws.Range("E19").Value = '  -0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.636.27'
$ws.Range("E20").Value = '  +0.85%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.646'
$ws.Range("E21").Value = '  +2.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.498'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.027'
$ws.Range("E23").Value = '  +1.69%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '142.89'
$ws.Range("E24").Value = '  -0.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.190'
$ws.Range("E25").Value = '  -6.85%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.693'
$ws.Range("E26").Value = '  +2.80%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.06'
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '111.79'
$ws.Range("E28").Value = '  +0.06%  '
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.126'
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08735'
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.04838'
$ws.Range("E32").Value = '  +0.59%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7230'
$ws.Range("E33").Value = '  +4.81%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.132'
$ws.Range("E34").Value = '  -0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.881'
$ws.Range("E35").Value = '  +0.83%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.117'
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.272'
$ws.Range("E37").Value = '  +2.41%  '
$ws.Range("E38").Value = '  +0.42%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.4878'
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9033'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '111.13'
$ws.Range("E41").Value = '  -1.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9996'
$ws.Range("E43").Value = '  -0.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.685'
$ws.Range("E44").Value = '  -1.28%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4217'
$ws.Range("E45").Value = '  +0.41%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.092'
$ws.Range("E46").Value = '  -0.45%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05888'
$ws.Range("E47").Value = '  +0.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1239'
$ws.Range("E48").Value = '  -2.23%  '
$ws.Range("E49").Value = '  -0.58%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.8892'
$ws.Range("E50").Value = '  +3.52%  '
$ws.Range("E51").Value = '  +1.53%  '
